$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet (fruit/vegetable weekly prices for Higo) has its weekly records
# re-permuted across rows 2-18, per the commit "Fruta / hortaliza, semanal".
# Column A-C, E-K stay fixed per row; columns D,L,M,N,O,P,Q,R,S,T are updated
# to the values that belonged (in the prior version) to a different row.

# Row 2 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 14
$ws.Range("D2").Value = 44320
$ws.Range('L2').Value = 'Primera'
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range('Q2').Value = '$/bandeja 7 kilos'
$ws.Range('R2').Value = 'Región Metropolitana'
$ws.Range("S2").Value = 1714
$ws.Range("T2").Value = 7

# Row 3 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 15
$ws.Range("D3").Value = 44320
$ws.Range('L3').Value = 'Segunda'
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range('Q3').Value = '$/bandeja 7 kilos'
$ws.Range('R3').Value = 'Región Metropolitana'
$ws.Range("S3").Value = 1143
$ws.Range("T3").Value = 7

# Row 4 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 16
$ws.Range("D4").Value = 44322
$ws.Range('L4').Value = 'Primera'
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range('Q4').Value = '$/bandeja 7 kilos'
$ws.Range('R4').Value = 'Región Metropolitana'
$ws.Range("S4").Value = 1714
$ws.Range("T4").Value = 7

# Row 5 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 17
$ws.Range("D5").Value = 44322
$ws.Range('L5').Value = 'Segunda'
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range('Q5').Value = '$/bandeja 7 kilos'
$ws.Range('R5').Value = 'Región Metropolitana'
$ws.Range("S5").Value = 1143
$ws.Range("T5").Value = 7

# Row 6 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 2
$ws.Range("D6").Value = 44980
$ws.Range('L6').Value = 'Primera'
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range('Q6').Value = '$/bandeja 7 kilos'
$ws.Range('R6').Value = 'Región Metropolitana'
$ws.Range("S6").Value = 2286
$ws.Range("T6").Value = 7

# Row 7 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 3
$ws.Range("D7").Value = 44980
$ws.Range('L7').Value = 'Segunda'
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 13000
$ws.Range('Q7').Value = '$/bandeja 7 kilos'
$ws.Range('R7').Value = 'Región Metropolitana'
$ws.Range("S7").Value = 1857
$ws.Range("T7").Value = 7

# Row 8 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 18
$ws.Range("D8").Value = 44971
$ws.Range('L8').Value = 'Primera'
$ws.Range("M8").Value = 25
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range('Q8').Value = '$/bandeja 5 kilos'
$ws.Range('R8').Value = 'Región Metropolitana'
$ws.Range("S8").Value = 3000
$ws.Range("T8").Value = 5

# Row 9 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 10
$ws.Range("D9").Value = 44299
$ws.Range('L9').Value = 'Primera'
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range('Q9').Value = '$/bandeja 7 kilos'
$ws.Range('R9').Value = 'Provincia de Santiago'
$ws.Range("S9").Value = 2143
$ws.Range("T9").Value = 7

# Row 10 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 11
$ws.Range("D10").Value = 44299
$ws.Range('L10').Value = 'Segunda'
$ws.Range("M10").Value = 75
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range('Q10').Value = '$/bandeja 7 kilos'
$ws.Range('R10').Value = 'Provincia de Santiago'
$ws.Range("S10").Value = 1714
$ws.Range("T10").Value = 7

# Row 11 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 8
$ws.Range("D11").Value = 44301
$ws.Range('L11').Value = 'Primera'
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 14000
$ws.Range('Q11').Value = '$/bandeja 7 kilos'
$ws.Range('R11').Value = 'Región Metropolitana'
$ws.Range("S11").Value = 2000
$ws.Range("T11").Value = 7

# Row 12 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 9
$ws.Range("D12").Value = 44301
$ws.Range('L12').Value = 'Segunda'
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range('Q12').Value = '$/bandeja 7 kilos'
$ws.Range('R12').Value = 'Región Metropolitana'
$ws.Range("S12").Value = 1714
$ws.Range("T12").Value = 7

# Row 13 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 4
$ws.Range("D13").Value = 44302
$ws.Range('L13').Value = 'Primera'
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range('Q13').Value = '$/bandeja 7 kilos'
$ws.Range('R13').Value = 'Región Metropolitana'
$ws.Range("S13").Value = 2143
$ws.Range("T13").Value = 7

# Row 14 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 5
$ws.Range("D14").Value = 44302
$ws.Range('L14').Value = 'Segunda'
$ws.Range("M14").Value = 30
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range('Q14').Value = '$/bandeja 7 kilos'
$ws.Range('R14').Value = 'Región Metropolitana'
$ws.Range("S14").Value = 1714
$ws.Range("T14").Value = 7

# Row 15 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 12
$ws.Range("D15").Value = 44292
$ws.Range('L15').Value = 'Primera'
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = 16000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 16000
$ws.Range('Q15').Value = '$/bandeja 7 kilos'
$ws.Range('R15').Value = 'Región Metropolitana'
$ws.Range("S15").Value = 2286
$ws.Range("T15").Value = 7

# Row 16 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 13
$ws.Range("D16").Value = 44292
$ws.Range('L16').Value = 'Segunda'
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range('Q16').Value = '$/bandeja 7 kilos'
$ws.Range('R16').Value = 'Región Metropolitana'
$ws.Range("S16").Value = 2143
$ws.Range("T16").Value = 7

# Row 17 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 6
$ws.Range("D17").Value = 44300
$ws.Range('L17').Value = 'Primera'
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 15000
$ws.Range('Q17').Value = '$/bandeja 7 kilos'
$ws.Range('R17').Value = 'Región Metropolitana'
$ws.Range("S17").Value = 2143
$ws.Range("T17").Value = 7

# Row 18 takes the D,L,M,N,O,P,Q,R,S,T values previously on row 7
$ws.Range("D18").Value = 44300
$ws.Range('L18').Value = 'Segunda'
$ws.Range("M18").Value = 80
$ws.Range("N18").Value = 12000
$ws.Range("O18").Value = 12000
$ws.Range("P18").Value = 12000
$ws.Range('Q18').Value = '$/bandeja 7 kilos'
$ws.Range('R18').Value = 'Región Metropolitana'
$ws.Range("S18").Value = 1714
$ws.Range("T18").Value = 7

